# Apply the edit described by the diff:
#  1. Remove the stray "_GoBack" bookmark from the very first paragraph
#     of the document (it keeps the "_gjdgxs" bookmark in place).
#  2. Insert a new paragraph right after the "A copy of any judgment
#     against me." paragraph containing "{% endif %}" (with the
#     "_GoBack" bookmark re-attached to it, collapsed at its start).
#  3. Remove the old, now-duplicate "{% endif %}" paragraph that used to
#     sit right before "Sincerely,".

$d = $word.ActiveDocument

function Find-ParaIndex($doc, [string]$matchText) {
    $paras = $doc.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $t = $paras.Item($i).Range.Text
        $tTrim = $t.TrimEnd([char]13, [char]7)
        if ($tTrim -eq $matchText) { return $i }
    }
    return -1
}

# --- 1. Drop the old _GoBack bookmark from paragraph 1 -----------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Insert the new "{% endif %}" paragraph -------------------------
$anchorIdx = Find-ParaIndex $d "A copy of any judgment against me."
if ($anchorIdx -eq -1) {
    throw "Could not locate 'A copy of any judgment against me.' paragraph"
}

# The paragraph right after the anchor (originally a blank paragraph).
$followingPara = $d.Paragraphs.Item($anchorIdx + 1)

# Insert a brand-new, unformatted paragraph before it -- this yields a
# plain <w:p><w:r/></w:p>, not inheriting the anchor paragraph's
# numbering/formatting.
$followingPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($anchorIdx + 1)
$newPara.Range.Text = "{% endif %}"

# Re-fetch after the text assignment (ranges/paragraph objects can be
# invalidated by the edit) and drop a collapsed "_GoBack" bookmark right
# at its start, matching the original bookmarkStart/bookmarkEnd pairing.
$newPara2 = $d.Paragraphs.Item($anchorIdx + 1)
$startPos = $newPara2.Range.Start
$bmRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 3. Remove the now-duplicate "{% endif %}" paragraph ---------------
# It is the standalone paragraph that used to sit directly before
# "Sincerely," (and right after the "...%}" run that closes the
# "{% if validation == False %}" block).
$sincerelyIdx = Find-ParaIndex $d "Sincerely,"
if ($sincerelyIdx -eq -1) {
    throw "Could not locate 'Sincerely,' paragraph"
}

$oldEndifPara = $d.Paragraphs.Item($sincerelyIdx - 1)
$oldEndifText = $oldEndifPara.Range.Text.TrimEnd([char]13, [char]7)
if ($oldEndifText -eq "{% endif %}") {
    $oldEndifPara.Range.Delete()
} else {
    throw "Unexpected paragraph before 'Sincerely,': $oldEndifText"
}
